$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "pos" (column A) and "size" (column D) values.
# Order of writes matters for shared-string table ordering: write A2 first,
# then D2:D10 (the new "size" value), then the remaining "pos" values so the
# resulting shared string table matches [-50,50], [30,30], [-30,50], ...
$ws.Range("A2").Value     = "[-50, 50]"
$ws.Range("D2:D10").Value = "[30, 30]"
$ws.Range("A3").Value     = "[-30, 50]"
$ws.Range("A4").Value     = "[-10, 50]"
$ws.Range("A5").Value     = "[-50, 30]"
$ws.Range("A6").Value     = "[-30, 30]"
$ws.Range("A7").Value     = "[-10, 30]"
$ws.Range("A8").Value     = "[-50, 10]"
$ws.Range("A9").Value     = "[-30, 10]"
$ws.Range("A10").Value    = "[-10, 10]"
